# Generate Report for Handback
#
# A new source file, c9c6d314-d009-4069-9d7e-65101a1f0b84.md, has been
# handed back in sync with en-US for both the zh-cn and de-de locales.
# Append one row describing it to each of the three report tables:
# Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileId   = "c9c6d314-d009-4069-9d7e-65101a1f0b84"
$fileName = "$fileId.md"
$pathName = "e2e\$fileName"
$status   = "Handed back: in sync with en-US"

$zhXlf = "$fileId.b3aa5cd8cf9eef0c031088cf60d20b21660fa35e.zh-cn.xlf"
$deXlf = "$fileId.b3aa5cd8cf9eef0c031088cf60d20b21660fa35e.de-de.xlf"

$overviewDate   = "2016-09-03 04:48:50"
$zhHandoffDate  = "2016-09-03 04:48:45"
$zhHandbackDate = "2016-09-03 04:49:07"
$deHandoffDate  = "2016-09-03 04:48:50"
$deHandbackDate = "2016-09-03 04:49:14"

$srcBaseUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master"
$zhcnBaseUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master"
$dedeBaseUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master"

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview  = $wb.Worksheets.Item("Overview")
$loOverview  = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rOverview   = $rowOverview.Range.Row

$wsOverview.Cells.Item($rOverview, 1).Value = $fileName
$wsOverview.Cells.Item($rOverview, 2).Value = $pathName
$wsOverview.Cells.Item($rOverview, 3).Value = ".md"
$wsOverview.Cells.Item($rOverview, 5).Value = $status
$wsOverview.Cells.Item($rOverview, 6).Value = $status
$wsOverview.Cells.Item($rOverview, 7).Value = $overviewDate
$wsOverview.Cells.Item($rOverview, 7).NumberFormat = $dateFormat

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($rOverview, 2), "$srcBaseUrl/e2e/$fileName", "", "", $pathName) | Out-Null

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh  = $wb.Worksheets.Item("zh-cn")
$loZh  = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()
$rZh   = $rowZh.Range.Row

$wsZh.Cells.Item($rZh, 1).Value = $fileName
$wsZh.Cells.Item($rZh, 2).Value = ".md"
$wsZh.Cells.Item($rZh, 3).Value = $status
$wsZh.Cells.Item($rZh, 4).Value = "e2e"
$wsZh.Cells.Item($rZh, 5).Value = "ht"
$wsZh.Cells.Item($rZh, 6).Value = "'True"
$wsZh.Cells.Item($rZh, 7).Value = $zhXlf
$wsZh.Cells.Item($rZh, 8).Value = $zhHandoffDate
$wsZh.Cells.Item($rZh, 8).NumberFormat = $dateFormat
$wsZh.Cells.Item($rZh, 9).Value = $fileName
$wsZh.Cells.Item($rZh, 10).Value = $zhXlf
$wsZh.Cells.Item($rZh, 11).Value = $zhHandbackDate
$wsZh.Cells.Item($rZh, 11).NumberFormat = $dateFormat
$wsZh.Cells.Item($rZh, 12).Value = "'"
$wsZh.Cells.Item($rZh, 13).Value = "'True"
$wsZh.Cells.Item($rZh, 14).Value = "'"
$wsZh.Cells.Item($rZh, 15).Value = "'False"
$wsZh.Cells.Item($rZh, 16).Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Cells.Item($rZh, 1), "$srcBaseUrl/e2e/$fileName", "", "", $fileName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Cells.Item($rZh, 9), "$zhcnBaseUrl/e2e/$fileName", "", "", $fileName) | Out-Null

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe  = $wb.Worksheets.Item("de-de")
$loDe  = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()
$rDe   = $rowDe.Range.Row

$wsDe.Cells.Item($rDe, 1).Value = $fileName
$wsDe.Cells.Item($rDe, 2).Value = ".md"
$wsDe.Cells.Item($rDe, 3).Value = $status
$wsDe.Cells.Item($rDe, 4).Value = "e2e"
$wsDe.Cells.Item($rDe, 5).Value = "ht"
$wsDe.Cells.Item($rDe, 6).Value = "'True"
$wsDe.Cells.Item($rDe, 7).Value = $deXlf
$wsDe.Cells.Item($rDe, 8).Value = $deHandoffDate
$wsDe.Cells.Item($rDe, 8).NumberFormat = $dateFormat
$wsDe.Cells.Item($rDe, 9).Value = $fileName
$wsDe.Cells.Item($rDe, 10).Value = $deXlf
$wsDe.Cells.Item($rDe, 11).Value = $deHandbackDate
$wsDe.Cells.Item($rDe, 11).NumberFormat = $dateFormat
$wsDe.Cells.Item($rDe, 12).Value = "'"
$wsDe.Cells.Item($rDe, 13).Value = "'True"
$wsDe.Cells.Item($rDe, 14).Value = "'"
$wsDe.Cells.Item($rDe, 15).Value = "'False"
$wsDe.Cells.Item($rDe, 16).Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Cells.Item($rDe, 1), "$srcBaseUrl/e2e/$fileName", "", "", $fileName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Cells.Item($rDe, 9), "$dedeBaseUrl/e2e/$fileName", "", "", $fileName) | Out-Null
